$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '61.832.91'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.425.14'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '405.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.85%  '
$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.415.31'
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.676'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.08%  '
$ws.Range("E11").Value = '  +16.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.141'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").Value = '3.980.63'
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.61'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("D17").Value = '3.460.90'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '61.902.90'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.79%  '
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000136'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +16.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '307.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '29.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.93%  '
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.80%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.174'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.64'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0481'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.95'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.125'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.283'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.24'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.60%  '
$ws.Range("D50").Value = '3.770.35'
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("D51").Value = '2.146.37'
$ws.Range("E51").Value = '  -1.24%  '
